$wb = $excel.ActiveWorkbook

# The "想去人数" (number of people interested) figures were refreshed for two
# events that appear on both the "展览" (Exhibition) sheet and the
# "全部类型" (All Types) sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 450
    $ws.Range("F3").Value = 5452
}
